# Hindalco: update PDFs & Excel (2025-08-12 08:31:40 UTC)
#
# Inserts 21 historical "Basic Price" rows (new rows 3-23) above the single
# existing data row, extending the sheet's used range from A1:F2 to A1:F23.
# Each new row mirrors the schema of row 2 (Sl.no., Description, Grade,
# Basic Price, Circular Date, Circular Link) but has no circular PDF link.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stamp formatting for the new rows before touching any values: copying
# row 2's cell formats down across A3:F23 gives every new cell the same
# styles as the existing data row (centered "General" style for
# A/B/C/E/F, the centered "0.000" numeric style for D) while leaving the
# cells themselves empty.
$ws.Range("A2:F2").Copy()
$ws.Range("A3:F23").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$desc = "2. P0610 (99.85% min) /P1020/ EC Grade Ingot & Sow 99.7% (min) / Cast Bar"

$ws.Range("A3").Value = 21
$ws.Range("B3").Value = $desc
$ws.Range("C3").Value = "P1020"
$ws.Range("D3").Value = 265.25
# Force the circular date into row 3 as literal text (matching the
# "dd.mm.yyyy" label style used throughout column E) instead of letting
# it auto-convert to a date serial number.
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "12.08.2025"
$ws.Range("E3").NumberFormat = "General"
# F3 intentionally left blank: no circular PDF link for this historical row.

$ws.Range("A4").Value = 22
$ws.Range("B4").Value = $desc
$ws.Range("C4").Value = "P1020"
$ws.Range("D4").Value = 268.5
# Force the circular date into row 4 as literal text (matching the
# "dd.mm.yyyy" label style used throughout column E) instead of letting
# it auto-convert to a date serial number.
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "08.08.2025"
$ws.Range("E4").NumberFormat = "General"
# F4 intentionally left blank: no circular PDF link for this historical row.

$ws.Range("A5").Value = 15
$ws.Range("B5").Value = $desc
$ws.Range("C5").Value = "P1020"
$ws.Range("D5").Value = 265.75
# Force the circular date into row 5 as literal text (matching the
# "dd.mm.yyyy" label style used throughout column E) instead of letting
# it auto-convert to a date serial number.
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "07.08.2025"
$ws.Range("E5").NumberFormat = "General"
# F5 intentionally left blank: no circular PDF link for this historical row.

$ws.Range("A6").Value = 20
$ws.Range("B6").Value = $desc
$ws.Range("C6").Value = "P1020"
$ws.Range("D6").Value = 263.75
# Force the circular date into row 6 as literal text (matching the
# "dd.mm.yyyy" label style used throughout column E) instead of letting
# it auto-convert to a date serial number.
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "05.08.2025"
$ws.Range("E6").NumberFormat = "General"
# F6 intentionally left blank: no circular PDF link for this historical row.

$ws.Range("A7").Value = 19
$ws.Range("B7").Value = $desc
$ws.Range("C7").Value = "P1020"
$ws.Range("D7").Value = 260.5
# Force the circular date into row 7 as literal text (matching the
# "dd.mm.yyyy" label style used throughout column E) instead of letting
# it auto-convert to a date serial number.
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "02.08.2025"
$ws.Range("E7").NumberFormat = "General"
# F7 intentionally left blank: no circular PDF link for this historical row.

$ws.Range("A8").Value = 18
$ws.Range("B8").Value = $desc
$ws.Range("C8").Value = "P1020"
$ws.Range("D8").Value = 264.5
# Force the circular date into row 8 as literal text (matching the
# "dd.mm.yyyy" label style used throughout column E) instead of letting
# it auto-convert to a date serial number.
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "01.08.2025"
$ws.Range("E8").NumberFormat = "General"
# F8 intentionally left blank: no circular PDF link for this historical row.

$ws.Range("A9").Value = 13
$ws.Range("B9").Value = $desc
$ws.Range("C9").Value = "P1020"
$ws.Range("D9").Value = 266.25
# Force the circular date into row 9 as literal text (matching the
# "dd.mm.yyyy" label style used throughout column E) instead of letting
# it auto-convert to a date serial number.
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "29.07.2025"
$ws.Range("E9").NumberFormat = "General"
# F9 intentionally left blank: no circular PDF link for this historical row.

$ws.Range("A10").Value = 17
$ws.Range("B10").Value = $desc
$ws.Range("C10").Value = "P1020"
$ws.Range("D10").Value = 268.5
# Force the circular date into row 10 as literal text (matching the
# "dd.mm.yyyy" label style used throughout column E) instead of letting
# it auto-convert to a date serial number.
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "26.07.2025"
$ws.Range("E10").NumberFormat = "General"
# F10 intentionally left blank: no circular PDF link for this historical row.

$ws.Range("A11").Value = 14
$ws.Range("B11").Value = $desc
$ws.Range("C11").Value = "P1020"
$ws.Range("D11").Value = 267
# Force the circular date into row 11 as literal text (matching the
# "dd.mm.yyyy" label style used throughout column E) instead of letting
# it auto-convert to a date serial number.
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "22.07.2025"
$ws.Range("E11").NumberFormat = "General"
# F11 intentionally left blank: no circular PDF link for this historical row.

$ws.Range("A12").Value = 16
$ws.Range("B12").Value = $desc
$ws.Range("C12").Value = "P1020"
$ws.Range("D12").Value = 261.5
# Force the circular date into row 12 as literal text (matching the
# "dd.mm.yyyy" label style used throughout column E) instead of letting
# it auto-convert to a date serial number.
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "19.07.2025"
$ws.Range("E12").NumberFormat = "General"
# F12 intentionally left blank: no circular PDF link for this historical row.

$ws.Range("A13").Value = 8
$ws.Range("B13").Value = $desc
$ws.Range("C13").Value = "P1020"
$ws.Range("D13").Value = 258
# Force the circular date into row 13 as literal text (matching the
# "dd.mm.yyyy" label style used throughout column E) instead of letting
# it auto-convert to a date serial number.
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "17.07.2025"
$ws.Range("E13").NumberFormat = "General"
# F13 intentionally left blank: no circular PDF link for this historical row.

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = $desc
$ws.Range("C14").Value = "P1020"
$ws.Range("D14").Value = 261.25
# Force the circular date into row 14 as literal text (matching the
# "dd.mm.yyyy" label style used throughout column E) instead of letting
# it auto-convert to a date serial number.
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "11.07.2025"
$ws.Range("E14").NumberFormat = "General"
# F14 intentionally left blank: no circular PDF link for this historical row.

$ws.Range("A15").Value = 10
$ws.Range("B15").Value = $desc
$ws.Range("C15").Value = "P1020"
$ws.Range("D15").Value = 258.5
# Force the circular date into row 15 as literal text (matching the
# "dd.mm.yyyy" label style used throughout column E) instead of letting
# it auto-convert to a date serial number.
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "05.07.2025"
$ws.Range("E15").NumberFormat = "General"
# F15 intentionally left blank: no circular PDF link for this historical row.

$ws.Range("A16").Value = 7
$ws.Range("B16").Value = $desc
$ws.Range("C16").Value = "P1020"
$ws.Range("D16").Value = 260.75
# Force the circular date into row 16 as literal text (matching the
# "dd.mm.yyyy" label style used throughout column E) instead of letting
# it auto-convert to a date serial number.
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "02.07.2025"
$ws.Range("E16").NumberFormat = "General"
# F16 intentionally left blank: no circular PDF link for this historical row.

$ws.Range("A17").Value = 9
$ws.Range("B17").Value = $desc
$ws.Range("C17").Value = "P1020"
$ws.Range("D17").Value = 263.25
# Force the circular date into row 17 as literal text (matching the
# "dd.mm.yyyy" label style used throughout column E) instead of letting
# it auto-convert to a date serial number.
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "28.06.2025"
$ws.Range("E17").NumberFormat = "General"
# F17 intentionally left blank: no circular PDF link for this historical row.

$ws.Range("A18").Value = 11
$ws.Range("B18").Value = $desc
$ws.Range("C18").Value = "P1020"
$ws.Range("D18").Value = 261.75
# Force the circular date into row 18 as literal text (matching the
# "dd.mm.yyyy" label style used throughout column E) instead of letting
# it auto-convert to a date serial number.
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "26.06.2025"
$ws.Range("E18").NumberFormat = "General"
# F18 intentionally left blank: no circular PDF link for this historical row.

$ws.Range("A19").Value = 6
$ws.Range("B19").Value = $desc
$ws.Range("C19").Value = "P1020"
$ws.Range("D19").Value = 264
# Force the circular date into row 19 as literal text (matching the
# "dd.mm.yyyy" label style used throughout column E) instead of letting
# it auto-convert to a date serial number.
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "25.06.2025"
$ws.Range("E19").NumberFormat = "General"
# F19 intentionally left blank: no circular PDF link for this historical row.

$ws.Range("A20").Value = 4
$ws.Range("B20").Value = $desc
$ws.Range("C20").Value = "P1020"
$ws.Range("D20").Value = 268.75
# Force the circular date into row 20 as literal text (matching the
# "dd.mm.yyyy" label style used throughout column E) instead of letting
# it auto-convert to a date serial number.
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "24.06.2025"
$ws.Range("E20").NumberFormat = "General"
# F20 intentionally left blank: no circular PDF link for this historical row.

$ws.Range("A21").Value = 2
$ws.Range("B21").Value = $desc
$ws.Range("C21").Value = "P1020"
$ws.Range("D21").Value = 262.25
# Force the circular date into row 21 as literal text (matching the
# "dd.mm.yyyy" label style used throughout column E) instead of letting
# it auto-convert to a date serial number.
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "19.06.2025"
$ws.Range("E21").NumberFormat = "General"
# F21 intentionally left blank: no circular PDF link for this historical row.

$ws.Range("A22").Value = 3
$ws.Range("B22").Value = $desc
$ws.Range("C22").Value = "P1020"
$ws.Range("D22").Value = 260
# Force the circular date into row 22 as literal text (matching the
# "dd.mm.yyyy" label style used throughout column E) instead of letting
# it auto-convert to a date serial number.
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "18.06.2025"
$ws.Range("E22").NumberFormat = "General"
# F22 intentionally left blank: no circular PDF link for this historical row.

$ws.Range("A23").Value = 5
$ws.Range("B23").Value = $desc
$ws.Range("C23").Value = "P1020"
$ws.Range("D23").Value = 256.5
# Force the circular date into row 23 as literal text (matching the
# "dd.mm.yyyy" label style used throughout column E) instead of letting
# it auto-convert to a date serial number.
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "12.06.2025"
$ws.Range("E23").NumberFormat = "General"
# F23 intentionally left blank: no circular PDF link for this historical row.
